$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 5
$ws.Range("P2").Value = 3
$ws.Range("Q2").ClearContents()

# Row 3
$ws.Range("M3").Value = 5
$ws.Range("P3").Value = 1
$ws.Range("Q3").ClearContents()

# Row 4
$ws.Range("M4").Value = 5
$ws.Range("P4").Value = 1
$ws.Range("Q4").ClearContents()

# Row 5
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 5
$ws.Range("P5").Value = 3
$ws.Range("Q5").ClearContents()
